# ---------------------------------------------------------------------------
# Applies the "Ajustes gitignore y cambios menores" edit to Data.xlsx:
#  - Reworks the "Inicio" sheet into a 3-column credentials table
#    (generarEvidencia / usuario / contrasenna -> SI / admin / 123456)
#  - Adds a new hidden sheet "op" with SI/NO options
#  - Adds a defined name "Producto" pointing at #REF!
#  - Restyles the header row (green fill) and the quoted "123456" value
#  - Re-themes the workbook from the "Aptos" Office theme to the classic
#    "Calibri / Office 2013-2022" theme colors & fonts
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Main "Inicio" sheet: new header/value table -----------------------
$ws1 = $wb.Worksheets.Item(1)
if ($ws1.Name -ne "Inicio") { $ws1.Name = "Inicio" }

$data = New-Object 'object[,]' 2,3
$data[0,0] = "generarEvidencia"
$data[0,1] = "usuario"
$data[0,2] = "contrasenna"
$data[1,0] = "SI"
$data[1,1] = "admin"
$data[1,2] = "'123456"
$ws1.Range("A1:C2").Value = $data

# Header row style: Calibri 10, green fill, vertically centered
$header = $ws1.Range("A1:C1")
$header.Interior.Color = 5296274
$header.Font.Name = "Calibri"
$header.Font.Size = 10
$header.VerticalAlignment = -4108

# Approximate the post-edit "best fit" column widths
$ws1.Columns.Item(1).ColumnWidth = 13.7369791666667
$ws1.Columns.Item(2).ColumnWidth = 17.1666666666667

$null = $ws1.Range("A3").Select()

# --- 2. New hidden "op" sheet ----------------------------------------------
$op = $wb.Worksheets.Add($null, $ws1)
$op.Name = "op"

$op.Range("A1").Value = "op"
$op.Range("A2").Value = "'SI"
$op.Range("A3").Value = "'NO"

$null = $op.Range("A2").Select()
$op.Visible = 0

# Re-activate Inicio so it stays the selected tab
$ws1.Activate()
$null = $ws1.Range("A3").Select()

# --- 3. Defined name: Producto -> #REF! ------------------------------------
$wb.Names.Add("Producto", "=#REF!")

# --- 4. Workbook view --------------------------------------------------
$excel.ActiveWindow.TabRatio = 0.788

# --- 5. Re-theme workbook (Aptos -> Calibri / Office 2013-2022 colors) ----
$theme = $wb.Theme
$colors = $theme.ThemeColorScheme
$colors.Colors(3).RGB  = 6968388    # dk2      44546A
$colors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Colors(5).RGB  = 12874308   # accent1  4472C4
$colors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407      # accent4  FFC000
$colors.Colors(9).RGB  = 13998939   # accent5  5B9BD5
$colors.Colors(10).RGB = 4697456    # accent6  70AD47
$colors.Colors(11).RGB = 12673797   # hlink    0563C1
$colors.Colors(12).RGB = 7491477    # folHlink 954F72

$fonts = $theme.ThemeFontScheme
$fonts.MajorFont.Item(1).Name = "Calibri Light"
$fonts.MinorFont.Item(1).Name = "Calibri"

Write-Host "edit applied"
